$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number (e.g. "518.32")
# must be forced to stay as TEXT (matching the inline-string cell type used
# throughout this sheet) instead of being auto-converted to a Number by Excel.
# We flip NumberFormat to Text ("@") just long enough to assign the value, then
# restore each cell's original Style so no formatting/style delta is introduced.
$textForceCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D15", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D48", "D50", "D51")
$origStyles = @{}
foreach ($addr in $textForceCells) { $origStyles[$addr] = $ws.Range($addr).Style }
foreach ($addr in $textForceCells) { $ws.Range($addr).NumberFormat = "@" }

# Apply the updated values (coin names / links / prices / volume changes).
$ws.Range("D2").Value = "57.295.41"
$ws.Range("E2").Value = "  -4.50%  "
$ws.Range("D3").Value = "3.135.14"
$ws.Range("E3").Value = "  -4.59%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "518.32"
$ws.Range("E5").Value = "  -6.67%  "
$ws.Range("D6").Value = "133.18"
$ws.Range("E6").Value = "  -6.04%  "
$ws.Range("D8").Value = "3.132.05"
$ws.Range("E8").Value = "  -4.65%  "
$ws.Range("D9").Value = "0.448"
$ws.Range("E9").Value = "  -5.84%  "
$ws.Range("D10").Value = "7.19"
$ws.Range("E10").Value = "  -7.65%  "
$ws.Range("D11").Value = "0.108"
$ws.Range("E11").Value = "  -8.45%  "
$ws.Range("D12").Value = "0.381"
$ws.Range("E12").Value = "  -5.61%  "
$ws.Range("D13").Value = "3.676.87"
$ws.Range("E13").Value = "  -4.61%  "
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").Value = "25.17"
$ws.Range("E15").Value = "  -6.74%  "
$ws.Range("D16").Value = "3.141.00"
$ws.Range("E16").Value = "  -4.56%  "
$ws.Range("D17").Value = "57.321.19"
$ws.Range("E17").Value = "  -4.56%  "
$ws.Range("D18").Value = "0.0000149"
$ws.Range("E18").Value = "  -9.45%  "
$ws.Range("D19").Value = "5.73"
$ws.Range("E19").Value = "  -5.80%  "
$ws.Range("D20").Value = "12.85"
$ws.Range("E20").Value = "  -9.16%  "
$ws.Range("D21").Value = "7.93"
$ws.Range("E21").Value = "  -6.21%  "
$ws.Range("D22").Value = "343.16"
$ws.Range("E22").Value = "  -7.60%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "68.72"
$ws.Range("E24").Value = "  -6.61%  "
$ws.Range("D25").Value = "0.502"
$ws.Range("E25").Value = "  -7.11%  "
$ws.Range("D26").Value = "3.271.24"
$ws.Range("E26").Value = "  -5.10%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "0.163"
$ws.Range("E28").Value = "  -5.49%  "
$ws.Range("D29").Value = "0.0₃0930"
$ws.Range("E29").Value = "  -9.00%  "
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").Value = "6.68"
$ws.Range("E31").Value = "  -5.87%  "
$ws.Range("D32").Value = "1.84"
$ws.Range("E32").Value = "  -8.77%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "21.60"
$ws.Range("E33").Value = "  -3.60%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "6.83"
$ws.Range("E34").Value = "  -10.15%  "
$ws.Range("D35").Value = "1.19"
$ws.Range("E35").Value = "  -3.69%  "
$ws.Range("D36").Value = "4.84"
$ws.Range("E36").Value = "  -6.40%  "
$ws.Range("D37").Value = "157.29"
$ws.Range("E37").Value = "  -5.54%  "
$ws.Range("D38").Value = "6.16"
$ws.Range("E38").Value = "  -7.71%  "
$ws.Range("D39").Value = "1.38"
$ws.Range("E39").Value = "  -8.26%  "
$ws.Range("D40").Value = "25.62"
$ws.Range("E40").Value = "  -4.39%  "
$ws.Range("D41").Value = "0.0683"
$ws.Range("E41").Value = "  -6.60%  "
$ws.Range("D42").Value = "3.166.13"
$ws.Range("E42").Value = "  -4.77%  "
$ws.Range("D43").Value = "40.31"
$ws.Range("E43").Value = "  -3.60%  "
$ws.Range("D44").Value = "0.693"
$ws.Range("E44").Value = "  -7.26%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "3.87"
$ws.Range("E45").Value = "  -6.88%  "
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").Value = "1.06"
$ws.Range("E46").Value = "  -5.30%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").Value = "1.45"
$ws.Range("E48").Value = "  -7.89%  "
$ws.Range("D49").Value = "2.234.33"
$ws.Range("E49").Value = "  -4.71%  "
$ws.Range("D50").Value = "6.13"
$ws.Range("E50").Value = "  -5.64%  "
$ws.Range("D51").Value = "19.93"
$ws.Range("E51").Value = "  -5.86%  "

# Restore original styles on the text-forced cells.
foreach ($addr in $textForceCells) { $ws.Range($addr).Style = $origStyles[$addr] }
